# Update TPM-derived metrics for Hgf-Sdc1 LR pairs (new TPM data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.001628901849080777
$ws.Range("J2").Value = 0.001628901849080777
$ws.Range("M2").Value = 3.173991666666667
$ws.Range("N2").Value = 9.521975000000001
$ws.Range("O2").Value = 0.1364420332266311
$ws.Range("P2").Value = 0.1364420332266311
$ws.Range("Q2").Value = 0.2972305656194445
$ws.Range("R2").Value = 2.675075090575
$ws.Range("S2").Value = 0.0002222506802152003
$ws.Range("T2").Value = 0.0002222506802152002
$ws.Range("I3").Value = 0.001628901849080777
$ws.Range("J3").Value = 0.001628901849080777
$ws.Range("O3").Value = 0.5095288789807429
$ws.Range("P3").Value = 0.5095288789807428
$ws.Range("S3").Value = 0.0008299725331317877
$ws.Range("T3").Value = 0.0008299725331317874
$ws.Range("I4").Value = 0.001628901849080777
$ws.Range("J4").Value = 0.001628901849080777
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4770453333333333
$ws.Range("N4").Value = 1.431136
$ws.Range("O4").Value = 0.02050699625485553
$ws.Range("P4").Value = 0.02050699625485553
$ws.Range("Q4").Value = 0.04467322827022222
$ws.Range("R4").Value = 0.402059054432
$ws.Range("S4").Value = 0.00003340388411862674
$ws.Range("T4").Value = 0.00003340388411862673
$ws.Range("I5").Value = 0.001628901849080777
$ws.Range("J5").Value = 0.001628901849080777
$ws.Range("M5").Value = 7.276137666666667
$ws.Range("N5").Value = 21.828413
$ws.Range("O5").Value = 0.3127831202907618
$ws.Range("P5").Value = 0.3127831202907618
$ws.Range("Q5").Value = 0.6813787625534445
$ws.Range("R5").Value = 6.132408862981
$ws.Range("S5").Value = 0.0005094930030028771
$ws.Range("T5").Value = 0.000509493003002877
$ws.Range("I6").Value = 0.001628901849080777
$ws.Range("J6").Value = 0.001628901849080777
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4824416666666667
$ws.Range("N6").Value = 1.447325
$ws.Range("O6").Value = 0.02073897124700851
$ws.Range("P6").Value = 0.02073897124700851
$ws.Range("Q6").Value = 0.04517857150277778
$ws.Range("R6").Value = 0.406607143525
$ws.Range("S6").Value = 0.00003378174861228523
$ws.Range("T6").Value = 0.00003378174861228523
$ws.Range("I7").Value = 0.1785947081647151
$ws.Range("J7").Value = 0.178594708164715
$ws.Range("M7").Value = 3.173991666666667
$ws.Range("N7").Value = 9.521975000000001
$ws.Range("O7").Value = 0.1364420332266311
$ws.Range("P7").Value = 0.1364420332266311
$ws.Range("Q7").Value = 32.58870763416112
$ws.Range("R7").Value = 293.29836870745
$ws.Range("S7").Value = 0.02436782510551055
$ws.Range("T7").Value = 0.02436782510551054
$ws.Range("I8").Value = 0.1785947081647151
$ws.Range("J8").Value = 0.178594708164715
$ws.Range("O8").Value = 0.5095288789807429
$ws.Range("P8").Value = 0.5095288789807428
$ws.Range("S8").Value = 0.09099916144306021
$ws.Range("T8").Value = 0.09099916144306017
$ws.Range("I9").Value = 0.1785947081647151
$ws.Range("J9").Value = 0.178594708164715
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4770453333333333
$ws.Range("N9").Value = 1.431136
$ws.Range("O9").Value = 0.02050699625485553
$ws.Range("P9").Value = 0.02050699625485553
$ws.Range("Q9").Value = 4.898025114403556
$ws.Range("R9").Value = 44.082226029632
$ws.Range("S9").Value = 0.003662441011470828
$ws.Range("T9").Value = 0.003662441011470827
$ws.Range("I10").Value = 0.1785947081647151
$ws.Range("J10").Value = 0.178594708164715
$ws.Range("M10").Value = 7.276137666666667
$ws.Range("N10").Value = 21.828413
$ws.Range("O10").Value = 0.3127831202907618
$ws.Range("P10").Value = 0.3127831202907618
$ws.Range("Q10").Value = 74.70716625224512
$ws.Range("R10").Value = 672.364496270206
$ws.Range("S10").Value = 0.05586141008717758
$ws.Range("T10").Value = 0.05586141008717757
$ws.Range("I11").Value = 0.1785947081647151
$ws.Range("J11").Value = 0.178594708164715
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.4824416666666667
$ws.Range("N11").Value = 1.447325
$ws.Range("O11").Value = 0.02073897124700851
$ws.Range("P11").Value = 0.02073897124700851
$ws.Range("Q11").Value = 4.953431538794444
$ws.Range("R11").Value = 44.58088384915
$ws.Range("S11").Value = 0.003703870517495902
$ws.Range("T11").Value = 0.003703870517495902
$ws.Range("G12").Value = 23.67539566666666
$ws.Range("H12").Value = 71.02618699999999
$ws.Range("I12").Value = 0.4118171950916292
$ws.Range("J12").Value = 0.4118171950916292
$ws.Range("M12").Value = 3.173991666666667
$ws.Range("N12").Value = 9.521975000000001
$ws.Range("O12").Value = 0.1364420332266311
$ws.Range("P12").Value = 0.1364420332266311
$ws.Range("Q12").Value = 75.14550855103612
$ws.Range("R12").Value = 676.309576959325
$ws.Range("S12").Value = 0.05618917541599012
$ws.Range("T12").Value = 0.0561891754159901
$ws.Range("G13").Value = 23.67539566666666
$ws.Range("H13").Value = 71.02618699999999
$ws.Range("I13").Value = 0.4118171950916292
$ws.Range("J13").Value = 0.4118171950916292
$ws.Range("O13").Value = 0.5095288789807429
$ws.Range("P13").Value = 0.5095288789807428
$ws.Range("Q13").Value = 280.6232495000224
$ws.Range("R13").Value = 2525.609245500202
$ws.Range("S13").Value = 0.2098327537600318
$ws.Range("T13").Value = 0.2098327537600317
$ws.Range("G14").Value = 23.67539566666666
$ws.Range("H14").Value = 71.02618699999999
$ws.Range("I14").Value = 0.4118171950916292
$ws.Range("J14").Value = 0.4118171950916292
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.4770453333333333
$ws.Range("N14").Value = 1.431136
$ws.Range("O14").Value = 0.02050699625485553
$ws.Range("P14").Value = 0.02050699625485553
$ws.Range("Q14").Value = 11.29423701760355
$ws.Range("R14").Value = 101.648133158432
$ws.Range("S14").Value = 0.00844513367742915
$ws.Range("T14").Value = 0.008445133677429146
$ws.Range("G15").Value = 23.67539566666666
$ws.Range("H15").Value = 71.02618699999999
$ws.Range("I15").Value = 0.4118171950916292
$ws.Range("J15").Value = 0.4118171950916292
$ws.Range("M15").Value = 7.276137666666667
$ws.Range("N15").Value = 21.828413
$ws.Range("O15").Value = 0.3127831202907618
$ws.Range("P15").Value = 0.3127831202907618
$ws.Range("Q15").Value = 172.2654381834701
$ws.Range("R15").Value = 1550.388943651231
$ws.Range("S15").Value = 0.1288094672701492
$ws.Range("T15").Value = 0.1288094672701492
$ws.Range("G16").Value = 23.67539566666666
$ws.Range("H16").Value = 71.02618699999999
$ws.Range("I16").Value = 0.4118171950916292
$ws.Range("J16").Value = 0.4118171950916292
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.4824416666666667
$ws.Range("N16").Value = 1.447325
$ws.Range("O16").Value = 0.02073897124700851
$ws.Range("P16").Value = 0.02073897124700851
$ws.Range("Q16").Value = 11.42199734441944
$ws.Range("R16").Value = 102.797976099775
$ws.Range("S16").Value = 0.008540664968028994
$ws.Range("T16").Value = 0.008540664968028992
$ws.Range("G17").Value = 0.3314846666666666
$ws.Range("H17").Value = 0.9944539999999999
$ws.Range("I17").Value = 0.005765947381177186
$ws.Range("J17").Value = 0.005765947381177185
$ws.Range("M17").Value = 3.173991666666667
$ws.Range("N17").Value = 9.521975000000001
$ws.Range("O17").Value = 0.1364420332266311
$ws.Range("P17").Value = 0.1364420332266311
$ws.Range("Q17").Value = 1.052129569627778
$ws.Range("R17").Value = 9.46916612665
$ws.Range("S17").Value = 0.0007867175841655844
$ws.Range("T17").Value = 0.0007867175841655841
$ws.Range("G18").Value = 0.3314846666666666
$ws.Range("H18").Value = 0.9944539999999999
$ws.Range("I18").Value = 0.005765947381177186
$ws.Range("J18").Value = 0.005765947381177185
$ws.Range("O18").Value = 0.5095288789807429
$ws.Range("P18").Value = 0.5095288789807428
$ws.Range("Q18").Value = 3.929070737787111
$ws.Range("R18").Value = 35.361636640084
$ws.Range("S18").Value = 0.002937916705393162
$ws.Range("T18").Value = 0.002937916705393161
$ws.Range("G19").Value = 0.3314846666666666
$ws.Range("H19").Value = 0.9944539999999999
$ws.Range("I19").Value = 0.005765947381177186
$ws.Range("J19").Value = 0.005765947381177185
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.4770453333333333
$ws.Range("N19").Value = 1.431136
$ws.Range("O19").Value = 0.02050699625485553
$ws.Range("P19").Value = 0.02050699625485553
$ws.Range("Q19").Value = 0.1581332133048889
$ws.Range("R19").Value = 1.423198919744
$ws.Range("S19").Value = 0.0001182422613514946
$ws.Range("T19").Value = 0.0001182422613514945
$ws.Range("G20").Value = 0.3314846666666666
$ws.Range("H20").Value = 0.9944539999999999
$ws.Range("I20").Value = 0.005765947381177186
$ws.Range("J20").Value = 0.005765947381177185
$ws.Range("M20").Value = 7.276137666666667
$ws.Range("N20").Value = 21.828413
$ws.Range("O20").Value = 0.3127831202907618
$ws.Range("P20").Value = 0.3127831202907618
$ws.Range("Q20").Value = 2.411928069055778
$ws.Range("R20").Value = 21.707352621502
$ws.Range("S20").Value = 0.001803491013316947
$ws.Range("T20").Value = 0.001803491013316947
$ws.Range("G21").Value = 0.3314846666666666
$ws.Range("H21").Value = 0.9944539999999999
$ws.Range("I21").Value = 0.005765947381177186
$ws.Range("J21").Value = 0.005765947381177185
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.4824416666666667
$ws.Range("N21").Value = 1.447325
$ws.Range("O21").Value = 0.02073897124700851
$ws.Range("P21").Value = 0.02073897124700851
$ws.Range("Q21").Value = 0.1599220150611111
$ws.Range("R21").Value = 1.43929813555
$ws.Range("S21").Value = 0.0001195798169499977
$ws.Range("T21").Value = 0.0001195798169499977
$ws.Range("G22").Value = 23.12211433333333
$ws.Range("H22").Value = 69.366343
$ws.Range("I22").Value = 0.4021932475133977
$ws.Range("J22").Value = 0.4021932475133977
$ws.Range("M22").Value = 3.173991666666667
$ws.Range("N22").Value = 9.521975000000001
$ws.Range("O22").Value = 0.1364420332266311
$ws.Range("P22").Value = 0.1364420332266311
$ws.Range("Q22").Value = 73.3893982097139
$ws.Range("R22").Value = 660.504583887425
$ws.Range("S22").Value = 0.0548760644407497
$ws.Range("T22").Value = 0.05487606444074968
$ws.Range("G23").Value = 23.12211433333333
$ws.Range("H23").Value = 69.366343
$ws.Range("I23").Value = 0.4021932475133977
$ws.Range("J23").Value = 0.4021932475133977
$ws.Range("O23").Value = 0.5095288789807429
$ws.Range("P23").Value = 0.5095288789807428
$ws.Range("Q23").Value = 274.0652342577976
$ws.Range("R23").Value = 2466.587108320178
$ws.Range("S23").Value = 0.204929074539126
$ws.Range("T23").Value = 0.204929074539126
$ws.Range("G24").Value = 23.12211433333333
$ws.Range("H24").Value = 69.366343
$ws.Range("I24").Value = 0.4021932475133977
$ws.Range("J24").Value = 0.4021932475133977
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.4770453333333333
$ws.Range("N24").Value = 1.431136
$ws.Range("O24").Value = 0.02050699625485553
$ws.Range("P24").Value = 0.02050699625485553
$ws.Range("Q24").Value = 11.03029673951644
$ws.Range("R24").Value = 99.272670655648
$ws.Range("S24").Value = 0.00824777542048543
$ws.Range("T24").Value = 0.008247775420485428
$ws.Range("G25").Value = 23.12211433333333
$ws.Range("H25").Value = 69.366343
$ws.Range("I25").Value = 0.4021932475133977
$ws.Range("J25").Value = 0.4021932475133977
$ws.Range("M25").Value = 7.276137666666667
$ws.Range("N25").Value = 21.828413
$ws.Range("O25").Value = 0.3127831202907618
$ws.Range("P25").Value = 0.3127831202907618
$ws.Range("Q25").Value = 168.2396870337399
$ws.Range("R25").Value = 1514.157183303659
$ws.Range("S25").Value = 0.1257992589171152
$ws.Range("T25").Value = 0.1257992589171152
$ws.Range("G26").Value = 23.12211433333333
$ws.Range("H26").Value = 69.366343
$ws.Range("I26").Value = 0.4021932475133977
$ws.Range("J26").Value = 0.4021932475133977
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.4824416666666667
$ws.Range("N26").Value = 1.447325
$ws.Range("O26").Value = 0.02073897124700851
$ws.Range("P26").Value = 0.02073897124700851
$ws.Range("Q26").Value = 11.15507137583056
$ws.Range("R26").Value = 100.395642382475
$ws.Range("S26").Value = 0.008341074195921333
$ws.Range("T26").Value = 0.008341074195921333
